$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Row 13: new diary entry for 1/19/2020 (serial 43849)
# ---------------------------------------------------------------------------
$ws.Range("A12:G12").Copy()
$ws.Range("A13:G13").PasteSpecial(-4122)   # xlPasteFormats - reuse existing styles
$excel.CutCopyMode = $false

$ws.Range("A13").Value = 43849
$ws.Range("B13").Value = "2:00 - 3:30 p.m."
$ws.Range("D13").Value = "Finish the pacman homework"
$ws.Range("E13").Value = "Finished the pacman hw"
$ws.Range("F13").Value = "PacMan is not a big app, but we still need to take some time reading the code. Fortunately we don't need to read every line of  code to make modification to the application"
$ws.Range("G13").Value = "Cool"
$ws.Rows("13").RowHeight = 85

# ---------------------------------------------------------------------------
# Row 14: new diary entry for 1/23/2020 (serial 43853)
# ---------------------------------------------------------------------------
$ws.Range("A12:G12").Copy()
$ws.Range("A14:G14").PasteSpecial(-4122)   # xlPasteFormats - reuse existing styles
$excel.CutCopyMode = $false

$ws.Range("A14").Value = 43853
$ws.Range("B14").Value = "5:00 - 7:50 p.m."
$ws.Range("D14").Value = "Learn about mental models, how to externalize mental model and  how to model the code using UML graph"
$ws.Range("G14").Value = "Worried, because the homework sounds intimidating"
$ws.Range("E14").Value = "Explored the pacman features by locating the code that implemented the feature. Built the UML graph of pacman using intellij plugin"
$ws.Range("F14").Value = "Mental model exist in our daily life. Learning how to externalize it is important for software engineers to commnuicate their ideas with other people"
$ws.Rows("14").RowHeight = 85

# ---------------------------------------------------------------------------
# Selection moves to B18, matching where the author left off editing
# ---------------------------------------------------------------------------
$ws.Range("B18").Select()
